$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells remain plain text (avoid Excel auto-numeric conversion)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Cells.Item(2, 4).Value = "25.912.01"
$ws.Cells.Item(2, 5).Value = "  +0.18%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "1.733.85"
$ws.Cells.Item(3, 5).Value = "  -0.34%  "

# Row 4
$ws.Cells.Item(4, 4).Value = "0.9994"
$ws.Cells.Item(4, 5).Value = "  -0.03%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "245.90"
$ws.Cells.Item(5, 5).Value = "  +3.19%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "1.000"
$ws.Cells.Item(6, 5).Value = "  +0.03%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "0.5033"
$ws.Cells.Item(7, 5).Value = "  -2.40%  "

# Row 8
$ws.Cells.Item(8, 4).Value = "0.2729"
$ws.Cells.Item(8, 5).Value = "  -0.25%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "0.06167"
$ws.Cells.Item(9, 5).Value = "  +0.61%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "1.739.97"
$ws.Cells.Item(10, 5).Value = "  +0.06%  "

# Row 11
$ws.Cells.Item(11, 4).Value = "0.07248"
$ws.Cells.Item(11, 5).Value = "  +0.90%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "0.6533"
$ws.Cells.Item(12, 5).Value = "  +1.26%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "15.17"
$ws.Cells.Item(13, 5).Value = "  +1.75%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "4.770"
$ws.Cells.Item(14, 5).Value = "  +3.90%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "77.06"
$ws.Cells.Item(15, 5).Value = "  -0.29%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "0.9986"
$ws.Cells.Item(16, 5).Value = "  -0.15%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "0.9989"
$ws.Cells.Item(17, 5).Value = "  -0.07%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "25.925.69"

# Row 19
$ws.Cells.Item(19, 5).Value = "  +1.25%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "0.000006811"
$ws.Cells.Item(20, 5).Value = "  +0.68%  "

# Row 21
$ws.Cells.Item(21, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(21, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(21, 4).Value = "1.963.85"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "

# Row 22
$ws.Cells.Item(22, 2).Value = "Uniswap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Cells.Item(22, 4).Value = "4.587"
$ws.Cells.Item(22, 5).Value = "  +7.58%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "8.787"
$ws.Cells.Item(23, 5).Value = "  +1.38%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "5.470"
$ws.Cells.Item(24, 5).Value = "  +4.40%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -3.55%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "15.25"
$ws.Cells.Item(26, 5).Value = "  +0.74%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "1.442"
$ws.Cells.Item(27, 5).Value = "  -4.44%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "1.788"
$ws.Cells.Item(28, 5).Value = "  +1.66%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "105.40"
$ws.Cells.Item(29, 5).Value = "  -0.42%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "3.989"
$ws.Cells.Item(30, 5).Value = "  +0.66%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "0.08115"
$ws.Cells.Item(31, 5).Value = "  -2.31%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "3.696"
$ws.Cells.Item(32, 5).Value = "  +1.47%  "

# Row 33
$ws.Cells.Item(33, 4).Value = "0.04730"
$ws.Cells.Item(33, 5).Value = "  +3.10%  "

# Row 34
$ws.Cells.Item(34, 4).Value = "2.654"
$ws.Cells.Item(34, 5).Value = "  -0.24%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "0.9968"
$ws.Cells.Item(35, 5).Value = "  +0.69%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "0.6119"
$ws.Cells.Item(36, 5).Value = "  -1.02%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "2.752"
$ws.Cells.Item(37, 5).Value = "  +2.31%  "

# Row 38
$ws.Cells.Item(38, 2).Value = "TrustWalletToken"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(38, 4).Value = "0.8859"
$ws.Cells.Item(38, 5).Value = "  +20.28%  "

# Row 39
$ws.Cells.Item(39, 2).Value = "VeChain"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(39, 4).Value = "0.01603"
$ws.Cells.Item(39, 5).Value = "  -0.48%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "1.954"
$ws.Cells.Item(40, 5).Value = "  +1.31%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "0.9991"
$ws.Cells.Item(41, 5).Value = "  -0.05%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "100.58"
$ws.Cells.Item(42, 5).Value = "  +2.97%  "

# Row 43
$ws.Cells.Item(43, 4).Value = "0.3908"
$ws.Cells.Item(43, 5).Value = "  +2.03%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "5.011"
$ws.Cells.Item(44, 5).Value = "  +1.18%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "0.1177"
$ws.Cells.Item(45, 5).Value = "  +4.67%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "6.320"
$ws.Cells.Item(46, 5).Value = "  +2.50%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "55.78"
$ws.Cells.Item(47, 5).Value = "  +1.79%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "0.05278"
$ws.Cells.Item(48, 5).Value = "  +0.37%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "30.76"
$ws.Cells.Item(49, 5).Value = "  +0.97%  "

# Row 50
$ws.Cells.Item(50, 2).Value = "EnergySwap"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(50, 4).Value = "7.689"
$ws.Cells.Item(50, 5).Value = "  +1.64%  "

# Row 51
$ws.Cells.Item(51, 2).Value = "Decentraland"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(51, 4).Value = "0.3485"
$ws.Cells.Item(51, 5).Value = "  +2.35%  "
